$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new formatted rows below the existing last row (row 7),
# copying row formatting from the row above (mirrors rows 2-7 style pattern).
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()

# Row 8 - batch_007
$ws.Range("A8").Value = "batch_007"
$ws.Range("B8").Value = "y"
$ws.Range("C8").Value = "批量操作语句7执行"
$ws.Range("D8").Value = "batchsql"
$ws.Range("F8").Value = "batch07"
$ws.Range("H8").Value = "batch_sql_07"
$ws.Range("I8").Value = "select * from `$batch07"
$ws.Range("J8").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_007.csv"
$ws.Range("M8").Value = "csv_containsAll"

# Row 9 - batch_008
$ws.Range("A9").Value = "batch_008"
$ws.Range("B9").Value = "y"
$ws.Range("C9").Value = "批量操作语句8执行"
$ws.Range("D9").Value = "batchsql"
$ws.Range("F9").Value = "batch08"
$ws.Range("H9").Value = "batch_sql_08"
$ws.Range("I9").Value = "select * from `$batch08"
$ws.Range("J9").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_008.csv"
$ws.Range("M9").Value = "csv_containsAll"

# Update selection / view (drops topLeftCell, moves selection to E7)
$ws.Range("E7").Select()
